$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: string cells are written in a specific order below so that the
# shared-string table indices line up with the target workbook (Excel
# assigns the next shared-string index the first time a distinct text
# value is written).

# --- Row 1 headers (A1 new, B1:E1 unchanged) ---
$ws.Range("A1").Value = "Drive Cycle"
$ws.Range("H1").Value = "BattDisEnrgyTotal ACC Wh"
$ws.Range("I1").Value = "BattRegEnrgyTotal ACC Wh"
$ws.Range("F1").Value = "BattDisEnrgyTotal NACC Wh"
$ws.Range("G1").Value = "BattRegEnrgyTotal NACC Wh"

# --- Row 2 (HWFET) updated values ---
$ws.Range("B2").Value = 204.89259999999999
$ws.Range("C2").Value = 204.7423
$ws.Range("E2").Value = 204.56829999999999

# --- Row 3 (UDDS) updated values ---
$ws.Range("B3").Value = 178.1069
$ws.Range("C3").Value = 176.6885
$ws.Range("E3").Value = 175.77449999999999

# --- Row 4 (US06) updated values ---
$ws.Range("B4").Value = 274.64490000000001
$ws.Range("C4").Value = 274.31979999999999
$ws.Range("E4").Value = 273.07990000000001

# --- New columns F:I for rows 2-4, formatted with 0.00 number format ---
$ws.Range("F2:I4").NumberFormat = "0.00"

$ws.Range("F2").Value = 2255.1999999999998
$ws.Range("G2").Value = 162.23099999999999
$ws.Range("H2").Value = 2232.6
$ws.Range("I2").Value = 140.1337

$ws.Range("F3").Value = 1717.8
$ws.Range("G3").Value = 401.03030000000001
$ws.Range("H3").Value = 1683
$ws.Range("I3").Value = 374.91039999999998

$ws.Range("F4").Value = 2762.6
$ws.Range("G4").Value = 575.33810000000005
$ws.Range("H4").Value = 2700.8
$ws.Range("I4").Value = 512.99760000000003

# --- Row 6: Non-ACC / ACC summary headers ---
$ws.Range("B6").Value = "BattDisEnrgyTotal "
$ws.Range("C6").Value = "BattRegEnrgyTotal"

# --- Row 7: Non ACC ---
$ws.Range("A7").Value = "Non ACC"
$ws.Range("B7").Value = 325.92989999999998
$ws.Range("C7").Value = 180.0763
$ws.Range("D7").Formula = "=B7-C7"

# --- Row 8: ACC ---
$ws.Range("A8").Value = "ACC"
$ws.Range("B8").Value = 314.33999999999997
$ws.Range("C8").Value = 176.43129999999999
$ws.Range("D8").Formula = "=B8-C8"

$ws.Range("D6").Value = "Total Energy Consumption"

# --- Row 10: WOT per-event headers ---
$ws.Range("A10").Value = "WOT"
$ws.Range("B10").Value = "BattDisEnrgy Per Event"
$ws.Range("C10").Value = "BattRegEnrgyTotal Per Event"
$ws.Range("D10").Value = "Total Energy Consumption Per event"

# --- Rows 11-15: per-event data ---
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 336.04390000000001
$ws.Range("C11").Value = 100.8359
$ws.Range("D11").Formula = "=B11-C11"

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 325.53609999999998
$ws.Range("C12").Value = 96.105500000000006
$ws.Range("D12").Formula = "=B12-C12"

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = 322.11500000000001
$ws.Range("C13").Value = 95.179000000000002
$ws.Range("D13").Formula = "=B13-C13"

$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 320.4556
$ws.Range("C14").Value = 95.211200000000005
$ws.Range("D14").Formula = "=B14-C14"

$ws.Range("A15").Value = 5
$ws.Range("B15").Value = 319.5206
$ws.Range("C15").Value = 95.635800000000003
$ws.Range("D15").Formula = "=B15-C15"

# --- Column widths for new columns F:I ---
$ws.Range("F1").ColumnWidth = 26
$ws.Range("G1").ColumnWidth = 26.5703125
$ws.Range("H1").ColumnWidth = 24.42578125
$ws.Range("I1").ColumnWidth = 25

# --- Selection matches the post-edit state ---
$ws.Range("E2:E4").Select()
